$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44957
$ws.Cells.Item(2, 9).Value = 'Primera'
$ws.Cells.Item(2, 10).Value = 30
$ws.Cells.Item(2, 11).Value = 45000
$ws.Cells.Item(2, 12).Value = 45000
$ws.Cells.Item(2, 13).Value = 45000
$ws.Cells.Item(2, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(2, 15).Value = 'Región del Maule'
$ws.Cells.Item(2, 16).Value = 1800
$ws.Cells.Item(2, 17).Value = 25

$ws.Cells.Item(3, 4).Value = 44214
$ws.Cells.Item(3, 9).Value = 'Primera'
$ws.Cells.Item(3, 10).Value = 70
$ws.Cells.Item(3, 11).Value = 35000
$ws.Cells.Item(3, 12).Value = 36000
$ws.Cells.Item(3, 13).Value = 35429
$ws.Cells.Item(3, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(3, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(3, 16).Value = 1417
$ws.Cells.Item(3, 17).Value = 25

$ws.Cells.Item(4, 4).Value = 44214
$ws.Cells.Item(4, 9).Value = 'Primera'
$ws.Cells.Item(4, 10).Value = 80
$ws.Cells.Item(4, 11).Value = 35000
$ws.Cells.Item(4, 12).Value = 35000
$ws.Cells.Item(4, 13).Value = 35000
$ws.Cells.Item(4, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(4, 15).Value = 'Región del Maule'
$ws.Cells.Item(4, 16).Value = 1400
$ws.Cells.Item(4, 17).Value = 25

$ws.Cells.Item(5, 4).Value = 44223
$ws.Cells.Item(5, 9).Value = 'Primera'
$ws.Cells.Item(5, 10).Value = 95
$ws.Cells.Item(5, 11).Value = 32000
$ws.Cells.Item(5, 12).Value = 32000
$ws.Cells.Item(5, 13).Value = 32000
$ws.Cells.Item(5, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(5, 15).Value = 'Región del Maule'
$ws.Cells.Item(5, 16).Value = 1280
$ws.Cells.Item(5, 17).Value = 25

$ws.Cells.Item(6, 4).Value = 44203
$ws.Cells.Item(6, 9).Value = 'Primera'
$ws.Cells.Item(6, 10).Value = 120
$ws.Cells.Item(6, 11).Value = 27000
$ws.Cells.Item(6, 12).Value = 27000
$ws.Cells.Item(6, 13).Value = 27000
$ws.Cells.Item(6, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(6, 15).Value = 'Región del Maule'
$ws.Cells.Item(6, 16).Value = 1080
$ws.Cells.Item(6, 17).Value = 25

$ws.Cells.Item(7, 4).Value = 44203
$ws.Cells.Item(7, 9).Value = 'Segunda'
$ws.Cells.Item(7, 10).Value = 40
$ws.Cells.Item(7, 11).Value = 25000
$ws.Cells.Item(7, 12).Value = 25000
$ws.Cells.Item(7, 13).Value = 25000
$ws.Cells.Item(7, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(7, 15).Value = 'Región del Maule'
$ws.Cells.Item(7, 16).Value = 1000
$ws.Cells.Item(7, 17).Value = 25

$ws.Cells.Item(8, 4).Value = 44225
$ws.Cells.Item(8, 9).Value = 'Primera'
$ws.Cells.Item(8, 10).Value = 115
$ws.Cells.Item(8, 11).Value = 28000
$ws.Cells.Item(8, 12).Value = 30000
$ws.Cells.Item(8, 13).Value = 29130
$ws.Cells.Item(8, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(8, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(8, 16).Value = 1165
$ws.Cells.Item(8, 17).Value = 25

$ws.Cells.Item(9, 4).Value = 44609
$ws.Cells.Item(9, 9).Value = 'Primera'
$ws.Cells.Item(9, 10).Value = 10
$ws.Cells.Item(9, 11).Value = 25000
$ws.Cells.Item(9, 12).Value = 25000
$ws.Cells.Item(9, 13).Value = 25000
$ws.Cells.Item(9, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(9, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(9, 16).Value = 1000
$ws.Cells.Item(9, 17).Value = 25

$ws.Cells.Item(10, 4).Value = 44585
$ws.Cells.Item(10, 9).Value = 'Primera'
$ws.Cells.Item(10, 10).Value = 25
$ws.Cells.Item(10, 11).Value = 28000
$ws.Cells.Item(10, 12).Value = 28000
$ws.Cells.Item(10, 13).Value = 28000
$ws.Cells.Item(10, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(10, 15).Value = 'Región del Maule'
$ws.Cells.Item(10, 16).Value = 1120
$ws.Cells.Item(10, 17).Value = 25

$ws.Cells.Item(11, 4).Value = 44602
$ws.Cells.Item(11, 9).Value = 'Primera'
$ws.Cells.Item(11, 10).Value = 110
$ws.Cells.Item(11, 11).Value = 25000
$ws.Cells.Item(11, 12).Value = 25000
$ws.Cells.Item(11, 13).Value = 25000
$ws.Cells.Item(11, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(11, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(11, 16).Value = 1000
$ws.Cells.Item(11, 17).Value = 25

$ws.Cells.Item(12, 4).Value = 44602
$ws.Cells.Item(12, 9).Value = 'Primera'
$ws.Cells.Item(12, 10).Value = 55
$ws.Cells.Item(12, 11).Value = 27000
$ws.Cells.Item(12, 12).Value = 27000
$ws.Cells.Item(12, 13).Value = 27000
$ws.Cells.Item(12, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(12, 15).Value = 'Región del Maule'
$ws.Cells.Item(12, 16).Value = 1080
$ws.Cells.Item(12, 17).Value = 25

$ws.Cells.Item(13, 4).Value = 44610
$ws.Cells.Item(13, 9).Value = 'Primera'
$ws.Cells.Item(13, 10).Value = 50
$ws.Cells.Item(13, 11).Value = 28000
$ws.Cells.Item(13, 12).Value = 28000
$ws.Cells.Item(13, 13).Value = 28000
$ws.Cells.Item(13, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(13, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(13, 16).Value = 1120
$ws.Cells.Item(13, 17).Value = 25

$ws.Cells.Item(14, 4).Value = 44991
$ws.Cells.Item(14, 9).Value = 'Primera'
$ws.Cells.Item(14, 10).Value = 65
$ws.Cells.Item(14, 11).Value = 38000
$ws.Cells.Item(14, 12).Value = 38000
$ws.Cells.Item(14, 13).Value = 38000
$ws.Cells.Item(14, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(14, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(14, 16).Value = 1520
$ws.Cells.Item(14, 17).Value = 25

$ws.Cells.Item(15, 4).Value = 45015
$ws.Cells.Item(15, 9).Value = 'Primera'
$ws.Cells.Item(15, 10).Value = 40
$ws.Cells.Item(15, 11).Value = 42000
$ws.Cells.Item(15, 12).Value = 42000
$ws.Cells.Item(15, 13).Value = 42000
$ws.Cells.Item(15, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(15, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(15, 16).Value = 1680
$ws.Cells.Item(15, 17).Value = 25

$ws.Cells.Item(16, 4).Value = 44579
$ws.Cells.Item(16, 9).Value = 'Primera'
$ws.Cells.Item(16, 10).Value = 50
$ws.Cells.Item(16, 11).Value = 28000
$ws.Cells.Item(16, 12).Value = 28000
$ws.Cells.Item(16, 13).Value = 28000
$ws.Cells.Item(16, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(16, 15).Value = 'Región del Maule'
$ws.Cells.Item(16, 16).Value = 1120
$ws.Cells.Item(16, 17).Value = 25

$ws.Cells.Item(17, 4).Value = 44959
$ws.Cells.Item(17, 9).Value = 'Primera'
$ws.Cells.Item(17, 10).Value = 50
$ws.Cells.Item(17, 11).Value = 40000
$ws.Cells.Item(17, 12).Value = 40000
$ws.Cells.Item(17, 13).Value = 40000
$ws.Cells.Item(17, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(17, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(17, 16).Value = 1600
$ws.Cells.Item(17, 17).Value = 25

$ws.Cells.Item(18, 4).Value = 44959
$ws.Cells.Item(18, 9).Value = 'Segunda'
$ws.Cells.Item(18, 10).Value = 20
$ws.Cells.Item(18, 11).Value = 35000
$ws.Cells.Item(18, 12).Value = 35000
$ws.Cells.Item(18, 13).Value = 35000
$ws.Cells.Item(18, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(18, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(18, 16).Value = 1400
$ws.Cells.Item(18, 17).Value = 25

$ws.Cells.Item(19, 4).Value = 44559
$ws.Cells.Item(19, 9).Value = 'Primera'
$ws.Cells.Item(19, 10).Value = 25
$ws.Cells.Item(19, 11).Value = 28000
$ws.Cells.Item(19, 12).Value = 28000
$ws.Cells.Item(19, 13).Value = 28000
$ws.Cells.Item(19, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(19, 15).Value = 'Región del Maule'
$ws.Cells.Item(19, 16).Value = 1120
$ws.Cells.Item(19, 17).Value = 25

$ws.Cells.Item(20, 4).Value = 44559
$ws.Cells.Item(20, 9).Value = 'Segunda'
$ws.Cells.Item(20, 10).Value = 30
$ws.Cells.Item(20, 11).Value = 25000
$ws.Cells.Item(20, 12).Value = 25000
$ws.Cells.Item(20, 13).Value = 25000
$ws.Cells.Item(20, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(20, 15).Value = 'Región del Maule'
$ws.Cells.Item(20, 16).Value = 1000
$ws.Cells.Item(20, 17).Value = 25

$ws.Cells.Item(21, 4).Value = 44196
$ws.Cells.Item(21, 9).Value = 'Segunda'
$ws.Cells.Item(21, 10).Value = 30
$ws.Cells.Item(21, 11).Value = 10000
$ws.Cells.Item(21, 12).Value = 10000
$ws.Cells.Item(21, 13).Value = 10000
$ws.Cells.Item(21, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(21, 15).Value = 'Región del Maule'
$ws.Cells.Item(21, 16).Value = 400
$ws.Cells.Item(21, 17).Value = 25

$ws.Cells.Item(22, 4).Value = 44217
$ws.Cells.Item(22, 9).Value = 'Primera'
$ws.Cells.Item(22, 10).Value = 200
$ws.Cells.Item(22, 11).Value = 32000
$ws.Cells.Item(22, 12).Value = 32000
$ws.Cells.Item(22, 13).Value = 32000
$ws.Cells.Item(22, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(22, 15).Value = 'Región del Maule'
$ws.Cells.Item(22, 16).Value = 1280
$ws.Cells.Item(22, 17).Value = 25

$ws.Cells.Item(23, 4).Value = 44921
$ws.Cells.Item(23, 9).Value = 'Primera'
$ws.Cells.Item(23, 10).Value = 155
$ws.Cells.Item(23, 11).Value = 50000
$ws.Cells.Item(23, 12).Value = 50000
$ws.Cells.Item(23, 13).Value = 50000
$ws.Cells.Item(23, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(23, 15).Value = 'Región del Maule'
$ws.Cells.Item(23, 16).Value = 2000
$ws.Cells.Item(23, 17).Value = 25

$ws.Cells.Item(24, 4).Value = 44922
$ws.Cells.Item(24, 9).Value = 'Primera'
$ws.Cells.Item(24, 10).Value = 200
$ws.Cells.Item(24, 11).Value = 43000
$ws.Cells.Item(24, 12).Value = 43000
$ws.Cells.Item(24, 13).Value = 43000
$ws.Cells.Item(24, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(24, 15).Value = 'Región del Maule'
$ws.Cells.Item(24, 16).Value = 1720
$ws.Cells.Item(24, 17).Value = 25

$ws.Cells.Item(25, 4).Value = 44574
$ws.Cells.Item(25, 9).Value = 'Primera'
$ws.Cells.Item(25, 10).Value = 210
$ws.Cells.Item(25, 11).Value = 28000
$ws.Cells.Item(25, 12).Value = 28000
$ws.Cells.Item(25, 13).Value = 28000
$ws.Cells.Item(25, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(25, 15).Value = 'Región del Maule'
$ws.Cells.Item(25, 16).Value = 1120
$ws.Cells.Item(25, 17).Value = 25

$ws.Cells.Item(26, 4).Value = 44574
$ws.Cells.Item(26, 9).Value = 'Segunda'
$ws.Cells.Item(26, 10).Value = 85
$ws.Cells.Item(26, 11).Value = 25000
$ws.Cells.Item(26, 12).Value = 25000
$ws.Cells.Item(26, 13).Value = 25000
$ws.Cells.Item(26, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(26, 15).Value = 'Región del Maule'
$ws.Cells.Item(26, 16).Value = 1000
$ws.Cells.Item(26, 17).Value = 25

$ws.Cells.Item(27, 4).Value = 44553
$ws.Cells.Item(27, 9).Value = 'Primera'
$ws.Cells.Item(27, 10).Value = 95
$ws.Cells.Item(27, 11).Value = 45000
$ws.Cells.Item(27, 12).Value = 45000
$ws.Cells.Item(27, 13).Value = 45000
$ws.Cells.Item(27, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(27, 15).Value = 'Región del Maule'
$ws.Cells.Item(27, 16).Value = 1800
$ws.Cells.Item(27, 17).Value = 25

$ws.Cells.Item(28, 4).Value = 44266
$ws.Cells.Item(28, 9).Value = 'Primera'
$ws.Cells.Item(28, 10).Value = 65
$ws.Cells.Item(28, 11).Value = 25000
$ws.Cells.Item(28, 12).Value = 25000
$ws.Cells.Item(28, 13).Value = 25000
$ws.Cells.Item(28, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(28, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(28, 16).Value = 1000
$ws.Cells.Item(28, 17).Value = 25

$ws.Cells.Item(29, 4).Value = 44987
$ws.Cells.Item(29, 9).Value = 'Primera'
$ws.Cells.Item(29, 10).Value = 40
$ws.Cells.Item(29, 11).Value = 32000
$ws.Cells.Item(29, 12).Value = 32000
$ws.Cells.Item(29, 13).Value = 32000
$ws.Cells.Item(29, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(29, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(29, 16).Value = 1280
$ws.Cells.Item(29, 17).Value = 25

$ws.Cells.Item(30, 4).Value = 44614
$ws.Cells.Item(30, 9).Value = 'Primera'
$ws.Cells.Item(30, 10).Value = 40
$ws.Cells.Item(30, 11).Value = 25000
$ws.Cells.Item(30, 12).Value = 25000
$ws.Cells.Item(30, 13).Value = 25000
$ws.Cells.Item(30, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(30, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(30, 16).Value = 1000
$ws.Cells.Item(30, 17).Value = 25

$ws.Cells.Item(31, 4).Value = 44244
$ws.Cells.Item(31, 9).Value = 'Primera'
$ws.Cells.Item(31, 10).Value = 80
$ws.Cells.Item(31, 11).Value = 29000
$ws.Cells.Item(31, 12).Value = 29000
$ws.Cells.Item(31, 13).Value = 29000
$ws.Cells.Item(31, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(31, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(31, 16).Value = 1160
$ws.Cells.Item(31, 17).Value = 25

$ws.Cells.Item(32, 4).Value = 44271
$ws.Cells.Item(32, 9).Value = 'Primera'
$ws.Cells.Item(32, 10).Value = 30
$ws.Cells.Item(32, 11).Value = 23000
$ws.Cells.Item(32, 12).Value = 23000
$ws.Cells.Item(32, 13).Value = 23000
$ws.Cells.Item(32, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(32, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(32, 16).Value = 920
$ws.Cells.Item(32, 17).Value = 25

$ws.Cells.Item(33, 4).Value = 44967
$ws.Cells.Item(33, 9).Value = 'Primera'
$ws.Cells.Item(33, 10).Value = 35
$ws.Cells.Item(33, 11).Value = 45000
$ws.Cells.Item(33, 12).Value = 45000
$ws.Cells.Item(33, 13).Value = 45000
$ws.Cells.Item(33, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(33, 15).Value = 'Región del Maule'
$ws.Cells.Item(33, 16).Value = 1800
$ws.Cells.Item(33, 17).Value = 25

$ws.Cells.Item(34, 4).Value = 44552
$ws.Cells.Item(34, 9).Value = 'Primera'
$ws.Cells.Item(34, 10).Value = 110
$ws.Cells.Item(34, 11).Value = 45000
$ws.Cells.Item(34, 12).Value = 45000
$ws.Cells.Item(34, 13).Value = 45000
$ws.Cells.Item(34, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(34, 15).Value = 'Región del Maule'
$ws.Cells.Item(34, 16).Value = 1800
$ws.Cells.Item(34, 17).Value = 25

$ws.Cells.Item(35, 4).Value = 44231
$ws.Cells.Item(35, 9).Value = 'Primera'
$ws.Cells.Item(35, 10).Value = 95
$ws.Cells.Item(35, 11).Value = 30000
$ws.Cells.Item(35, 12).Value = 30000
$ws.Cells.Item(35, 13).Value = 30000
$ws.Cells.Item(35, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(35, 15).Value = 'Región del Maule'
$ws.Cells.Item(35, 16).Value = 1200
$ws.Cells.Item(35, 17).Value = 25

$ws.Cells.Item(36, 4).Value = 44208
$ws.Cells.Item(36, 9).Value = 'Primera'
$ws.Cells.Item(36, 10).Value = 95
$ws.Cells.Item(36, 11).Value = 27000
$ws.Cells.Item(36, 12).Value = 28000
$ws.Cells.Item(36, 13).Value = 27526
$ws.Cells.Item(36, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(36, 15).Value = 'Región del Maule'
$ws.Cells.Item(36, 16).Value = 1101
$ws.Cells.Item(36, 17).Value = 25

$ws.Cells.Item(37, 4).Value = 44595
$ws.Cells.Item(37, 9).Value = 'Primera'
$ws.Cells.Item(37, 10).Value = 100
$ws.Cells.Item(37, 11).Value = 28000
$ws.Cells.Item(37, 12).Value = 28000
$ws.Cells.Item(37, 13).Value = 28000
$ws.Cells.Item(37, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(37, 15).Value = 'Región del Maule'
$ws.Cells.Item(37, 16).Value = 1120
$ws.Cells.Item(37, 17).Value = 25

$ws.Cells.Item(38, 4).Value = 44586
$ws.Cells.Item(38, 9).Value = 'Primera'
$ws.Cells.Item(38, 10).Value = 55
$ws.Cells.Item(38, 11).Value = 28000
$ws.Cells.Item(38, 12).Value = 28000
$ws.Cells.Item(38, 13).Value = 28000
$ws.Cells.Item(38, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(38, 15).Value = 'Región del Maule'
$ws.Cells.Item(38, 16).Value = 1120
$ws.Cells.Item(38, 17).Value = 25

$ws.Cells.Item(39, 4).Value = 44239
$ws.Cells.Item(39, 9).Value = 'Primera'
$ws.Cells.Item(39, 10).Value = 210
$ws.Cells.Item(39, 11).Value = 29000
$ws.Cells.Item(39, 12).Value = 29000
$ws.Cells.Item(39, 13).Value = 29000
$ws.Cells.Item(39, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(39, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(39, 16).Value = 1160
$ws.Cells.Item(39, 17).Value = 25

$ws.Cells.Item(40, 4).Value = 44235
$ws.Cells.Item(40, 9).Value = 'Primera'
$ws.Cells.Item(40, 10).Value = 210
$ws.Cells.Item(40, 11).Value = 29000
$ws.Cells.Item(40, 12).Value = 29000
$ws.Cells.Item(40, 13).Value = 29000
$ws.Cells.Item(40, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(40, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(40, 16).Value = 1160
$ws.Cells.Item(40, 17).Value = 25

$ws.Cells.Item(41, 4).Value = 44952
$ws.Cells.Item(41, 9).Value = 'Primera'
$ws.Cells.Item(41, 10).Value = 135
$ws.Cells.Item(41, 11).Value = 45000
$ws.Cells.Item(41, 12).Value = 45000
$ws.Cells.Item(41, 13).Value = 45000
$ws.Cells.Item(41, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(41, 15).Value = 'Región del Maule'
$ws.Cells.Item(41, 16).Value = 1800
$ws.Cells.Item(41, 17).Value = 25

$ws.Cells.Item(42, 4).Value = 44221
$ws.Cells.Item(42, 9).Value = 'Primera'
$ws.Cells.Item(42, 10).Value = 310
$ws.Cells.Item(42, 11).Value = 32000
$ws.Cells.Item(42, 12).Value = 34000
$ws.Cells.Item(42, 13).Value = 32806
$ws.Cells.Item(42, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(42, 15).Value = 'Región del Maule'
$ws.Cells.Item(42, 16).Value = 1312
$ws.Cells.Item(42, 17).Value = 25

$ws.Cells.Item(43, 4).Value = 44216
$ws.Cells.Item(43, 9).Value = 'Primera'
$ws.Cells.Item(43, 10).Value = 70
$ws.Cells.Item(43, 11).Value = 35000
$ws.Cells.Item(43, 12).Value = 35000
$ws.Cells.Item(43, 13).Value = 35000
$ws.Cells.Item(43, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(43, 15).Value = 'Región del Maule'
$ws.Cells.Item(43, 16).Value = 1400
$ws.Cells.Item(43, 17).Value = 25

$ws.Cells.Item(44, 4).Value = 44966
$ws.Cells.Item(44, 9).Value = 'Primera'
$ws.Cells.Item(44, 10).Value = 80
$ws.Cells.Item(44, 11).Value = 45000
$ws.Cells.Item(44, 12).Value = 45000
$ws.Cells.Item(44, 13).Value = 45000
$ws.Cells.Item(44, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(44, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(44, 16).Value = 1800
$ws.Cells.Item(44, 17).Value = 25

$ws.Cells.Item(45, 4).Value = 44966
$ws.Cells.Item(45, 9).Value = 'Primera'
$ws.Cells.Item(45, 10).Value = 30
$ws.Cells.Item(45, 11).Value = 45000
$ws.Cells.Item(45, 12).Value = 45000
$ws.Cells.Item(45, 13).Value = 45000
$ws.Cells.Item(45, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(45, 15).Value = 'Región del Maule'
$ws.Cells.Item(45, 16).Value = 1800
$ws.Cells.Item(45, 17).Value = 25

$ws.Cells.Item(46, 4).Value = 44966
$ws.Cells.Item(46, 9).Value = 'Segunda'
$ws.Cells.Item(46, 10).Value = 50
$ws.Cells.Item(46, 11).Value = 40000
$ws.Cells.Item(46, 12).Value = 40000
$ws.Cells.Item(46, 13).Value = 40000
$ws.Cells.Item(46, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(46, 15).Value = 'Región del Maule'
$ws.Cells.Item(46, 16).Value = 1600
$ws.Cells.Item(46, 17).Value = 25

$ws.Cells.Item(47, 4).Value = 44224
$ws.Cells.Item(47, 9).Value = 'Primera'
$ws.Cells.Item(47, 10).Value = 380
$ws.Cells.Item(47, 11).Value = 27000
$ws.Cells.Item(47, 12).Value = 27000
$ws.Cells.Item(47, 13).Value = 27000
$ws.Cells.Item(47, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(47, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(47, 16).Value = 1080
$ws.Cells.Item(47, 17).Value = 25

$ws.Cells.Item(48, 4).Value = 44224
$ws.Cells.Item(48, 9).Value = 'Primera'
$ws.Cells.Item(48, 10).Value = 330
$ws.Cells.Item(48, 11).Value = 28000
$ws.Cells.Item(48, 12).Value = 30000
$ws.Cells.Item(48, 13).Value = 28909
$ws.Cells.Item(48, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(48, 15).Value = 'Región del Maule'
$ws.Cells.Item(48, 16).Value = 1156
$ws.Cells.Item(48, 17).Value = 25

$ws.Cells.Item(49, 4).Value = 44277
$ws.Cells.Item(49, 9).Value = 'Primera'
$ws.Cells.Item(49, 10).Value = 95
$ws.Cells.Item(49, 11).Value = 23000
$ws.Cells.Item(49, 12).Value = 23000
$ws.Cells.Item(49, 13).Value = 23000
$ws.Cells.Item(49, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(49, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(49, 16).Value = 920
$ws.Cells.Item(49, 17).Value = 25

$ws.Cells.Item(50, 4).Value = 44195
$ws.Cells.Item(50, 9).Value = 'Primera'
$ws.Cells.Item(50, 10).Value = 155
$ws.Cells.Item(50, 11).Value = 25000
$ws.Cells.Item(50, 12).Value = 25000
$ws.Cells.Item(50, 13).Value = 25000
$ws.Cells.Item(50, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(50, 15).Value = 'Región del Maule'
$ws.Cells.Item(50, 16).Value = 1000
$ws.Cells.Item(50, 17).Value = 25

$ws.Cells.Item(51, 4).Value = 44914
$ws.Cells.Item(51, 9).Value = 'Primera'
$ws.Cells.Item(51, 10).Value = 110
$ws.Cells.Item(51, 11).Value = 50000
$ws.Cells.Item(51, 12).Value = 50000
$ws.Cells.Item(51, 13).Value = 50000
$ws.Cells.Item(51, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(51, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(51, 16).Value = 2000
$ws.Cells.Item(51, 17).Value = 25

$ws.Cells.Item(52, 4).Value = 44581
$ws.Cells.Item(52, 9).Value = 'Primera'
$ws.Cells.Item(52, 10).Value = 80
$ws.Cells.Item(52, 11).Value = 28000
$ws.Cells.Item(52, 12).Value = 29000
$ws.Cells.Item(52, 13).Value = 28625
$ws.Cells.Item(52, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(52, 15).Value = 'Región del Maule'
$ws.Cells.Item(52, 16).Value = 1145
$ws.Cells.Item(52, 17).Value = 25

$ws.Cells.Item(53, 4).Value = 44588
$ws.Cells.Item(53, 9).Value = 'Primera'
$ws.Cells.Item(53, 10).Value = 185
$ws.Cells.Item(53, 11).Value = 28000
$ws.Cells.Item(53, 12).Value = 28000
$ws.Cells.Item(53, 13).Value = 28000
$ws.Cells.Item(53, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(53, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(53, 16).Value = 1120
$ws.Cells.Item(53, 17).Value = 25

$ws.Cells.Item(54, 4).Value = 44189
$ws.Cells.Item(54, 9).Value = 'Primera'
$ws.Cells.Item(54, 10).Value = 10
$ws.Cells.Item(54, 11).Value = 30000
$ws.Cells.Item(54, 12).Value = 30000
$ws.Cells.Item(54, 13).Value = 30000
$ws.Cells.Item(54, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(54, 15).Value = 'Región del Maule'
$ws.Cells.Item(54, 16).Value = 1200
$ws.Cells.Item(54, 17).Value = 25

$ws.Cells.Item(55, 4).Value = 44259
$ws.Cells.Item(55, 9).Value = 'Primera'
$ws.Cells.Item(55, 10).Value = 100
$ws.Cells.Item(55, 11).Value = 25000
$ws.Cells.Item(55, 12).Value = 25000
$ws.Cells.Item(55, 13).Value = 25000
$ws.Cells.Item(55, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(55, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(55, 16).Value = 1000
$ws.Cells.Item(55, 17).Value = 25

$ws.Cells.Item(56, 4).Value = 44166
$ws.Cells.Item(56, 9).Value = 'Primera'
$ws.Cells.Item(56, 10).Value = 45
$ws.Cells.Item(56, 11).Value = 1700
$ws.Cells.Item(56, 12).Value = 1800
$ws.Cells.Item(56, 13).Value = 1756
$ws.Cells.Item(56, 14).Value = '$/kilo'
$ws.Cells.Item(56, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(56, 16).Value = 1756
$ws.Cells.Item(56, 17).Value = 1

$ws.Cells.Item(57, 4).Value = 44923
$ws.Cells.Item(57, 9).Value = 'Primera'
$ws.Cells.Item(57, 10).Value = 35
$ws.Cells.Item(57, 11).Value = 43000
$ws.Cells.Item(57, 12).Value = 43000
$ws.Cells.Item(57, 13).Value = 43000
$ws.Cells.Item(57, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(57, 15).Value = 'Región del Maule'
$ws.Cells.Item(57, 16).Value = 1720
$ws.Cells.Item(57, 17).Value = 25

$ws.Cells.Item(58, 4).Value = 44923
$ws.Cells.Item(58, 9).Value = 'Primera'
$ws.Cells.Item(58, 10).Value = 35
$ws.Cells.Item(58, 11).Value = 43000
$ws.Cells.Item(58, 12).Value = 43000
$ws.Cells.Item(58, 13).Value = 43000
$ws.Cells.Item(58, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(58, 15).Value = 'Región del Maule'
$ws.Cells.Item(58, 16).Value = 1720
$ws.Cells.Item(58, 17).Value = 25

$ws.Cells.Item(59, 4).Value = 44945
$ws.Cells.Item(59, 9).Value = 'Primera'
$ws.Cells.Item(59, 10).Value = 75
$ws.Cells.Item(59, 11).Value = 45000
$ws.Cells.Item(59, 12).Value = 45000
$ws.Cells.Item(59, 13).Value = 45000
$ws.Cells.Item(59, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(59, 15).Value = 'Región del Maule'
$ws.Cells.Item(59, 16).Value = 1800
$ws.Cells.Item(59, 17).Value = 25

$ws.Cells.Item(60, 4).Value = 44193
$ws.Cells.Item(60, 9).Value = 'Primera'
$ws.Cells.Item(60, 10).Value = 120
$ws.Cells.Item(60, 11).Value = 29000
$ws.Cells.Item(60, 12).Value = 30000
$ws.Cells.Item(60, 13).Value = 29542
$ws.Cells.Item(60, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(60, 15).Value = 'Región del Maule'
$ws.Cells.Item(60, 16).Value = 1182
$ws.Cells.Item(60, 17).Value = 25

$ws.Cells.Item(61, 4).Value = 44568
$ws.Cells.Item(61, 9).Value = 'Primera'
$ws.Cells.Item(61, 10).Value = 80
$ws.Cells.Item(61, 11).Value = 30000
$ws.Cells.Item(61, 12).Value = 30000
$ws.Cells.Item(61, 13).Value = 30000
$ws.Cells.Item(61, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(61, 15).Value = 'Región del Maule'
$ws.Cells.Item(61, 16).Value = 1200
$ws.Cells.Item(61, 17).Value = 25

$ws.Cells.Item(62, 4).Value = 44209
$ws.Cells.Item(62, 9).Value = 'Primera'
$ws.Cells.Item(62, 10).Value = 100
$ws.Cells.Item(62, 11).Value = 27000
$ws.Cells.Item(62, 12).Value = 27000
$ws.Cells.Item(62, 13).Value = 27000
$ws.Cells.Item(62, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(62, 15).Value = 'Región del Maule'
$ws.Cells.Item(62, 16).Value = 1080
$ws.Cells.Item(62, 17).Value = 25

$ws.Cells.Item(63, 4).Value = 44606
$ws.Cells.Item(63, 9).Value = 'Primera'
$ws.Cells.Item(63, 10).Value = 50
$ws.Cells.Item(63, 11).Value = 25000
$ws.Cells.Item(63, 12).Value = 25000
$ws.Cells.Item(63, 13).Value = 25000
$ws.Cells.Item(63, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(63, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(63, 16).Value = 1000
$ws.Cells.Item(63, 17).Value = 25

$ws.Cells.Item(64, 4).Value = 44236
$ws.Cells.Item(64, 9).Value = 'Primera'
$ws.Cells.Item(64, 10).Value = 110
$ws.Cells.Item(64, 11).Value = 29000
$ws.Cells.Item(64, 12).Value = 29000
$ws.Cells.Item(64, 13).Value = 29000
$ws.Cells.Item(64, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(64, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(64, 16).Value = 1160
$ws.Cells.Item(64, 17).Value = 25

$ws.Cells.Item(65, 4).Value = 44624
$ws.Cells.Item(65, 9).Value = 'Primera'
$ws.Cells.Item(65, 10).Value = 30
$ws.Cells.Item(65, 11).Value = 25000
$ws.Cells.Item(65, 12).Value = 25000
$ws.Cells.Item(65, 13).Value = 25000
$ws.Cells.Item(65, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(65, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(65, 16).Value = 1000
$ws.Cells.Item(65, 17).Value = 25

$ws.Cells.Item(66, 4).Value = 45043
$ws.Cells.Item(66, 9).Value = 'Primera'
$ws.Cells.Item(66, 10).Value = 30
$ws.Cells.Item(66, 11).Value = 35000
$ws.Cells.Item(66, 12).Value = 35000
$ws.Cells.Item(66, 13).Value = 35000
$ws.Cells.Item(66, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(66, 15).Value = 'Región del Maule'
$ws.Cells.Item(66, 16).Value = 1400
$ws.Cells.Item(66, 17).Value = 25

$ws.Cells.Item(67, 4).Value = 44587
$ws.Cells.Item(67, 9).Value = 'Primera'
$ws.Cells.Item(67, 10).Value = 20
$ws.Cells.Item(67, 11).Value = 28000
$ws.Cells.Item(67, 12).Value = 28000
$ws.Cells.Item(67, 13).Value = 28000
$ws.Cells.Item(67, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(67, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(67, 16).Value = 1120
$ws.Cells.Item(67, 17).Value = 25

$ws.Cells.Item(68, 4).Value = 44600
$ws.Cells.Item(68, 9).Value = 'Primera'
$ws.Cells.Item(68, 10).Value = 145
$ws.Cells.Item(68, 11).Value = 25000
$ws.Cells.Item(68, 12).Value = 27000
$ws.Cells.Item(68, 13).Value = 25828
$ws.Cells.Item(68, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(68, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(68, 16).Value = 1033
$ws.Cells.Item(68, 17).Value = 25

$ws.Cells.Item(69, 4).Value = 44589
$ws.Cells.Item(69, 9).Value = 'Primera'
$ws.Cells.Item(69, 10).Value = 140
$ws.Cells.Item(69, 11).Value = 25000
$ws.Cells.Item(69, 12).Value = 28000
$ws.Cells.Item(69, 13).Value = 26179
$ws.Cells.Item(69, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(69, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(69, 16).Value = 1047
$ws.Cells.Item(69, 17).Value = 25

$ws.Cells.Item(70, 4).Value = 44243
$ws.Cells.Item(70, 9).Value = 'Primera'
$ws.Cells.Item(70, 10).Value = 145
$ws.Cells.Item(70, 11).Value = 28000
$ws.Cells.Item(70, 12).Value = 29000
$ws.Cells.Item(70, 13).Value = 28448
$ws.Cells.Item(70, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(70, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(70, 16).Value = 1138
$ws.Cells.Item(70, 17).Value = 25

$ws.Cells.Item(71, 4).Value = 44944
$ws.Cells.Item(71, 9).Value = 'Primera'
$ws.Cells.Item(71, 10).Value = 35
$ws.Cells.Item(71, 11).Value = 45000
$ws.Cells.Item(71, 12).Value = 45000
$ws.Cells.Item(71, 13).Value = 45000
$ws.Cells.Item(71, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(71, 15).Value = 'Región del Maule'
$ws.Cells.Item(71, 16).Value = 1800
$ws.Cells.Item(71, 17).Value = 25

$ws.Cells.Item(72, 4).Value = 44204
$ws.Cells.Item(72, 9).Value = 'Primera'
$ws.Cells.Item(72, 10).Value = 50
$ws.Cells.Item(72, 11).Value = 27000
$ws.Cells.Item(72, 12).Value = 27000
$ws.Cells.Item(72, 13).Value = 27000
$ws.Cells.Item(72, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(72, 15).Value = 'Región del Maule'
$ws.Cells.Item(72, 16).Value = 1080
$ws.Cells.Item(72, 17).Value = 25

$ws.Cells.Item(73, 4).Value = 44194
$ws.Cells.Item(73, 9).Value = 'Primera'
$ws.Cells.Item(73, 10).Value = 350
$ws.Cells.Item(73, 11).Value = 22000
$ws.Cells.Item(73, 12).Value = 22000
$ws.Cells.Item(73, 13).Value = 22000
$ws.Cells.Item(73, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(73, 15).Value = 'Región del Maule'
$ws.Cells.Item(73, 16).Value = 880
$ws.Cells.Item(73, 17).Value = 25

$ws.Cells.Item(74, 4).Value = 44215
$ws.Cells.Item(74, 9).Value = 'Primera'
$ws.Cells.Item(74, 10).Value = 80
$ws.Cells.Item(74, 11).Value = 35000
$ws.Cells.Item(74, 12).Value = 35000
$ws.Cells.Item(74, 13).Value = 35000
$ws.Cells.Item(74, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(74, 15).Value = 'Región del Maule'
$ws.Cells.Item(74, 16).Value = 1400
$ws.Cells.Item(74, 17).Value = 25

$ws.Cells.Item(75, 4).Value = 44960
$ws.Cells.Item(75, 9).Value = 'Primera'
$ws.Cells.Item(75, 10).Value = 30
$ws.Cells.Item(75, 11).Value = 45000
$ws.Cells.Item(75, 12).Value = 45000
$ws.Cells.Item(75, 13).Value = 45000
$ws.Cells.Item(75, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(75, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(75, 16).Value = 1800
$ws.Cells.Item(75, 17).Value = 25

$ws.Cells.Item(76, 4).Value = 44578
$ws.Cells.Item(76, 9).Value = 'Primera'
$ws.Cells.Item(76, 10).Value = 100
$ws.Cells.Item(76, 11).Value = 28000
$ws.Cells.Item(76, 12).Value = 28000
$ws.Cells.Item(76, 13).Value = 28000
$ws.Cells.Item(76, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(76, 15).Value = 'Región del Maule'
$ws.Cells.Item(76, 16).Value = 1120
$ws.Cells.Item(76, 17).Value = 25

$ws.Cells.Item(77, 4).Value = 44951
$ws.Cells.Item(77, 9).Value = 'Primera'
$ws.Cells.Item(77, 10).Value = 55
$ws.Cells.Item(77, 11).Value = 45000
$ws.Cells.Item(77, 12).Value = 45000
$ws.Cells.Item(77, 13).Value = 45000
$ws.Cells.Item(77, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(77, 15).Value = 'Región del Maule'
$ws.Cells.Item(77, 16).Value = 1800
$ws.Cells.Item(77, 17).Value = 25

$ws.Cells.Item(78, 4).Value = 44970
$ws.Cells.Item(78, 9).Value = 'Primera'
$ws.Cells.Item(78, 10).Value = 140
$ws.Cells.Item(78, 11).Value = 30000
$ws.Cells.Item(78, 12).Value = 35000
$ws.Cells.Item(78, 13).Value = 32857
$ws.Cells.Item(78, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(78, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(78, 16).Value = 1314
$ws.Cells.Item(78, 17).Value = 25

$ws.Cells.Item(79, 4).Value = 44279
$ws.Cells.Item(79, 9).Value = 'Primera'
$ws.Cells.Item(79, 10).Value = 50
$ws.Cells.Item(79, 11).Value = 23000
$ws.Cells.Item(79, 12).Value = 23000
$ws.Cells.Item(79, 13).Value = 23000
$ws.Cells.Item(79, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(79, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(79, 16).Value = 920
$ws.Cells.Item(79, 17).Value = 25

$ws.Cells.Item(80, 4).Value = 44567
$ws.Cells.Item(80, 9).Value = 'Primera'
$ws.Cells.Item(80, 10).Value = 130
$ws.Cells.Item(80, 11).Value = 27000
$ws.Cells.Item(80, 12).Value = 28000
$ws.Cells.Item(80, 13).Value = 27615
$ws.Cells.Item(80, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(80, 15).Value = 'Región del Maule'
$ws.Cells.Item(80, 16).Value = 1105
$ws.Cells.Item(80, 17).Value = 25

$ws.Cells.Item(81, 4).Value = 44936
$ws.Cells.Item(81, 9).Value = 'Primera'
$ws.Cells.Item(81, 10).Value = 30
$ws.Cells.Item(81, 11).Value = 50000
$ws.Cells.Item(81, 12).Value = 50000
$ws.Cells.Item(81, 13).Value = 50000
$ws.Cells.Item(81, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(81, 15).Value = 'Región del Maule'
$ws.Cells.Item(81, 16).Value = 2000
$ws.Cells.Item(81, 17).Value = 25

$ws.Cells.Item(82, 4).Value = 44278
$ws.Cells.Item(82, 9).Value = 'Primera'
$ws.Cells.Item(82, 10).Value = 65
$ws.Cells.Item(82, 11).Value = 23000
$ws.Cells.Item(82, 12).Value = 23000
$ws.Cells.Item(82, 13).Value = 23000
$ws.Cells.Item(82, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(82, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(82, 16).Value = 920
$ws.Cells.Item(82, 17).Value = 25

$ws.Cells.Item(83, 4).Value = 44971
$ws.Cells.Item(83, 9).Value = 'Primera'
$ws.Cells.Item(83, 10).Value = 20
$ws.Cells.Item(83, 11).Value = 35000
$ws.Cells.Item(83, 12).Value = 35000
$ws.Cells.Item(83, 13).Value = 35000
$ws.Cells.Item(83, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(83, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(83, 16).Value = 1400
$ws.Cells.Item(83, 17).Value = 25

$ws.Cells.Item(84, 4).Value = 44267
$ws.Cells.Item(84, 9).Value = 'Primera'
$ws.Cells.Item(84, 10).Value = 65
$ws.Cells.Item(84, 11).Value = 25000
$ws.Cells.Item(84, 12).Value = 25000
$ws.Cells.Item(84, 13).Value = 25000
$ws.Cells.Item(84, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(84, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(84, 16).Value = 1000
$ws.Cells.Item(84, 17).Value = 25

$ws.Cells.Item(85, 4).Value = 44603
$ws.Cells.Item(85, 9).Value = 'Primera'
$ws.Cells.Item(85, 10).Value = 65
$ws.Cells.Item(85, 11).Value = 25000
$ws.Cells.Item(85, 12).Value = 25000
$ws.Cells.Item(85, 13).Value = 25000
$ws.Cells.Item(85, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(85, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(85, 16).Value = 1000
$ws.Cells.Item(85, 17).Value = 25

$ws.Cells.Item(86, 4).Value = 44974
$ws.Cells.Item(86, 9).Value = 'Primera'
$ws.Cells.Item(86, 10).Value = 30
$ws.Cells.Item(86, 11).Value = 30000
$ws.Cells.Item(86, 12).Value = 30000
$ws.Cells.Item(86, 13).Value = 30000
$ws.Cells.Item(86, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(86, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(86, 16).Value = 1200
$ws.Cells.Item(86, 17).Value = 25

$ws.Cells.Item(87, 4).Value = 44188
$ws.Cells.Item(87, 9).Value = 'Primera'
$ws.Cells.Item(87, 10).Value = 50
$ws.Cells.Item(87, 11).Value = 35000
$ws.Cells.Item(87, 12).Value = 35000
$ws.Cells.Item(87, 13).Value = 35000
$ws.Cells.Item(87, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(87, 15).Value = 'Región del Maule'
$ws.Cells.Item(87, 16).Value = 1400
$ws.Cells.Item(87, 17).Value = 25

$ws.Cells.Item(88, 4).Value = 44943
$ws.Cells.Item(88, 9).Value = 'Primera'
$ws.Cells.Item(88, 10).Value = 65
$ws.Cells.Item(88, 11).Value = 45000
$ws.Cells.Item(88, 12).Value = 45000
$ws.Cells.Item(88, 13).Value = 45000
$ws.Cells.Item(88, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(88, 15).Value = 'Región del Maule'
$ws.Cells.Item(88, 16).Value = 1800
$ws.Cells.Item(88, 17).Value = 25

$ws.Cells.Item(89, 4).Value = 44637
$ws.Cells.Item(89, 9).Value = 'Primera'
$ws.Cells.Item(89, 10).Value = 50
$ws.Cells.Item(89, 11).Value = 20000
$ws.Cells.Item(89, 12).Value = 20000
$ws.Cells.Item(89, 13).Value = 20000
$ws.Cells.Item(89, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(89, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(89, 16).Value = 800
$ws.Cells.Item(89, 17).Value = 25

$ws.Cells.Item(90, 4).Value = 44935
$ws.Cells.Item(90, 9).Value = 'Primera'
$ws.Cells.Item(90, 10).Value = 55
$ws.Cells.Item(90, 11).Value = 50000
$ws.Cells.Item(90, 12).Value = 50000
$ws.Cells.Item(90, 13).Value = 50000
$ws.Cells.Item(90, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(90, 15).Value = 'Región del Maule'
$ws.Cells.Item(90, 16).Value = 2000
$ws.Cells.Item(90, 17).Value = 25

$ws.Cells.Item(91, 4).Value = 44200
$ws.Cells.Item(91, 9).Value = 'Primera'
$ws.Cells.Item(91, 10).Value = 80
$ws.Cells.Item(91, 11).Value = 32000
$ws.Cells.Item(91, 12).Value = 32000
$ws.Cells.Item(91, 13).Value = 32000
$ws.Cells.Item(91, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(91, 15).Value = 'Región del Maule'
$ws.Cells.Item(91, 16).Value = 1280
$ws.Cells.Item(91, 17).Value = 25

$ws.Cells.Item(92, 4).Value = 44949
$ws.Cells.Item(92, 9).Value = 'Primera'
$ws.Cells.Item(92, 10).Value = 75
$ws.Cells.Item(92, 11).Value = 45000
$ws.Cells.Item(92, 12).Value = 45000
$ws.Cells.Item(92, 13).Value = 45000
$ws.Cells.Item(92, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(92, 15).Value = 'Región del Maule'
$ws.Cells.Item(92, 16).Value = 1800
$ws.Cells.Item(92, 17).Value = 25

$ws.Cells.Item(93, 4).Value = 44249
$ws.Cells.Item(93, 9).Value = 'Primera'
$ws.Cells.Item(93, 10).Value = 175
$ws.Cells.Item(93, 11).Value = 28000
$ws.Cells.Item(93, 12).Value = 29000
$ws.Cells.Item(93, 13).Value = 28543
$ws.Cells.Item(93, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(93, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(93, 16).Value = 1142
$ws.Cells.Item(93, 17).Value = 25

$ws.Cells.Item(94, 4).Value = 44937
$ws.Cells.Item(94, 9).Value = 'Primera'
$ws.Cells.Item(94, 10).Value = 65
$ws.Cells.Item(94, 11).Value = 45000
$ws.Cells.Item(94, 12).Value = 45000
$ws.Cells.Item(94, 13).Value = 45000
$ws.Cells.Item(94, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(94, 15).Value = 'Región del Maule'
$ws.Cells.Item(94, 16).Value = 1800
$ws.Cells.Item(94, 17).Value = 25

$ws.Cells.Item(95, 4).Value = 44937
$ws.Cells.Item(95, 9).Value = 'Segunda'
$ws.Cells.Item(95, 10).Value = 55
$ws.Cells.Item(95, 11).Value = 40000
$ws.Cells.Item(95, 12).Value = 40000
$ws.Cells.Item(95, 13).Value = 40000
$ws.Cells.Item(95, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(95, 15).Value = 'Región del Maule'
$ws.Cells.Item(95, 16).Value = 1600
$ws.Cells.Item(95, 17).Value = 25

$ws.Cells.Item(96, 4).Value = 44202
$ws.Cells.Item(96, 9).Value = 'Primera'
$ws.Cells.Item(96, 10).Value = 40
$ws.Cells.Item(96, 11).Value = 28000
$ws.Cells.Item(96, 12).Value = 28000
$ws.Cells.Item(96, 13).Value = 28000
$ws.Cells.Item(96, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(96, 15).Value = 'Región del Maule'
$ws.Cells.Item(96, 16).Value = 1120
$ws.Cells.Item(96, 17).Value = 25

$ws.Cells.Item(97, 4).Value = 44938
$ws.Cells.Item(97, 9).Value = 'Primera'
$ws.Cells.Item(97, 10).Value = 35
$ws.Cells.Item(97, 11).Value = 50000
$ws.Cells.Item(97, 12).Value = 50000
$ws.Cells.Item(97, 13).Value = 50000
$ws.Cells.Item(97, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(97, 15).Value = 'Región del Maule'
$ws.Cells.Item(97, 16).Value = 2000
$ws.Cells.Item(97, 17).Value = 25

$ws.Cells.Item(98, 4).Value = 44218
$ws.Cells.Item(98, 9).Value = 'Primera'
$ws.Cells.Item(98, 10).Value = 80
$ws.Cells.Item(98, 11).Value = 35000
$ws.Cells.Item(98, 12).Value = 35000
$ws.Cells.Item(98, 13).Value = 35000
$ws.Cells.Item(98, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(98, 15).Value = 'Región del Maule'
$ws.Cells.Item(98, 16).Value = 1400
$ws.Cells.Item(98, 17).Value = 25

$ws.Cells.Item(99, 4).Value = 44273
$ws.Cells.Item(99, 9).Value = 'Primera'
$ws.Cells.Item(99, 10).Value = 40
$ws.Cells.Item(99, 11).Value = 25000
$ws.Cells.Item(99, 12).Value = 25000
$ws.Cells.Item(99, 13).Value = 25000
$ws.Cells.Item(99, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(99, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(99, 16).Value = 1000
$ws.Cells.Item(99, 17).Value = 25

$ws.Cells.Item(100, 4).Value = 44616
$ws.Cells.Item(100, 9).Value = 'Primera'
$ws.Cells.Item(100, 10).Value = 90
$ws.Cells.Item(100, 11).Value = 23000
$ws.Cells.Item(100, 12).Value = 25000
$ws.Cells.Item(100, 13).Value = 23889
$ws.Cells.Item(100, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(100, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(100, 16).Value = 956
$ws.Cells.Item(100, 17).Value = 25

$ws.Cells.Item(101, 4).Value = 44638
$ws.Cells.Item(101, 9).Value = 'Primera'
$ws.Cells.Item(101, 10).Value = 30
$ws.Cells.Item(101, 11).Value = 20000
$ws.Cells.Item(101, 12).Value = 20000
$ws.Cells.Item(101, 13).Value = 20000
$ws.Cells.Item(101, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(101, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(101, 16).Value = 800
$ws.Cells.Item(101, 17).Value = 25

$ws.Cells.Item(102, 4).Value = 44596
$ws.Cells.Item(102, 9).Value = 'Primera'
$ws.Cells.Item(102, 10).Value = 40
$ws.Cells.Item(102, 11).Value = 28000
$ws.Cells.Item(102, 12).Value = 28000
$ws.Cells.Item(102, 13).Value = 28000
$ws.Cells.Item(102, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(102, 15).Value = 'Región del Maule'
$ws.Cells.Item(102, 16).Value = 1120
$ws.Cells.Item(102, 17).Value = 25

$ws.Cells.Item(103, 4).Value = 44942
$ws.Cells.Item(103, 9).Value = 'Primera'
$ws.Cells.Item(103, 10).Value = 45
$ws.Cells.Item(103, 11).Value = 50000
$ws.Cells.Item(103, 12).Value = 50000
$ws.Cells.Item(103, 13).Value = 50000
$ws.Cells.Item(103, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(103, 15).Value = 'Región del Maule'
$ws.Cells.Item(103, 16).Value = 2000
$ws.Cells.Item(103, 17).Value = 25

$ws.Cells.Item(104, 4).Value = 44964
$ws.Cells.Item(104, 9).Value = 'Primera'
$ws.Cells.Item(104, 10).Value = 50
$ws.Cells.Item(104, 11).Value = 45000
$ws.Cells.Item(104, 12).Value = 45000
$ws.Cells.Item(104, 13).Value = 45000
$ws.Cells.Item(104, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(104, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(104, 16).Value = 1800
$ws.Cells.Item(104, 17).Value = 25

$ws.Cells.Item(105, 4).Value = 44964
$ws.Cells.Item(105, 9).Value = 'Primera'
$ws.Cells.Item(105, 10).Value = 30
$ws.Cells.Item(105, 11).Value = 45000
$ws.Cells.Item(105, 12).Value = 45000
$ws.Cells.Item(105, 13).Value = 45000
$ws.Cells.Item(105, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(105, 15).Value = 'Región del Maule'
$ws.Cells.Item(105, 16).Value = 1800
$ws.Cells.Item(105, 17).Value = 25

$ws.Cells.Item(106, 4).Value = 44222
$ws.Cells.Item(106, 9).Value = 'Primera'
$ws.Cells.Item(106, 10).Value = 110
$ws.Cells.Item(106, 11).Value = 32000
$ws.Cells.Item(106, 12).Value = 32000
$ws.Cells.Item(106, 13).Value = 32000
$ws.Cells.Item(106, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(106, 15).Value = 'Región del Maule'
$ws.Cells.Item(106, 16).Value = 1280
$ws.Cells.Item(106, 17).Value = 25

$ws.Cells.Item(107, 4).Value = 44582
$ws.Cells.Item(107, 9).Value = 'Primera'
$ws.Cells.Item(107, 10).Value = 50
$ws.Cells.Item(107, 11).Value = 28000
$ws.Cells.Item(107, 12).Value = 28000
$ws.Cells.Item(107, 13).Value = 28000
$ws.Cells.Item(107, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(107, 15).Value = 'Región del Maule'
$ws.Cells.Item(107, 16).Value = 1120
$ws.Cells.Item(107, 17).Value = 25

$ws.Cells.Item(108, 4).Value = 44956
$ws.Cells.Item(108, 9).Value = 'Primera'
$ws.Cells.Item(108, 10).Value = 50
$ws.Cells.Item(108, 11).Value = 50000
$ws.Cells.Item(108, 12).Value = 50000
$ws.Cells.Item(108, 13).Value = 50000
$ws.Cells.Item(108, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(108, 15).Value = 'Región del Maule'
$ws.Cells.Item(108, 16).Value = 2000
$ws.Cells.Item(108, 17).Value = 25

$ws.Cells.Item(109, 4).Value = 44242
$ws.Cells.Item(109, 9).Value = 'Primera'
$ws.Cells.Item(109, 10).Value = 155
$ws.Cells.Item(109, 11).Value = 29000
$ws.Cells.Item(109, 12).Value = 29000
$ws.Cells.Item(109, 13).Value = 29000
$ws.Cells.Item(109, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(109, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(109, 16).Value = 1160
$ws.Cells.Item(109, 17).Value = 25

$ws.Cells.Item(110, 4).Value = 44601
$ws.Cells.Item(110, 9).Value = 'Primera'
$ws.Cells.Item(110, 10).Value = 65
$ws.Cells.Item(110, 11).Value = 25000
$ws.Cells.Item(110, 12).Value = 25000
$ws.Cells.Item(110, 13).Value = 25000
$ws.Cells.Item(110, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(110, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(110, 16).Value = 1000
$ws.Cells.Item(110, 17).Value = 25

$ws.Cells.Item(111, 4).Value = 44575
$ws.Cells.Item(111, 9).Value = 'Primera'
$ws.Cells.Item(111, 10).Value = 125
$ws.Cells.Item(111, 11).Value = 28000
$ws.Cells.Item(111, 12).Value = 28000
$ws.Cells.Item(111, 13).Value = 28000
$ws.Cells.Item(111, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(111, 15).Value = 'Región del Maule'
$ws.Cells.Item(111, 16).Value = 1120
$ws.Cells.Item(111, 17).Value = 25

$ws.Cells.Item(112, 4).Value = 44918
$ws.Cells.Item(112, 9).Value = 'Primera'
$ws.Cells.Item(112, 10).Value = 65
$ws.Cells.Item(112, 11).Value = 5000
$ws.Cells.Item(112, 12).Value = 5000
$ws.Cells.Item(112, 13).Value = 5000
$ws.Cells.Item(112, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(112, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(112, 16).Value = 200
$ws.Cells.Item(112, 17).Value = 25

$ws.Cells.Item(113, 4).Value = 44623
$ws.Cells.Item(113, 9).Value = 'Primera'
$ws.Cells.Item(113, 10).Value = 90
$ws.Cells.Item(113, 11).Value = 25000
$ws.Cells.Item(113, 12).Value = 25000
$ws.Cells.Item(113, 13).Value = 25000
$ws.Cells.Item(113, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(113, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(113, 16).Value = 1000
$ws.Cells.Item(113, 17).Value = 25

$ws.Cells.Item(114, 4).Value = 44238
$ws.Cells.Item(114, 9).Value = 'Primera'
$ws.Cells.Item(114, 10).Value = 220
$ws.Cells.Item(114, 11).Value = 28000
$ws.Cells.Item(114, 12).Value = 29000
$ws.Cells.Item(114, 13).Value = 28568
$ws.Cells.Item(114, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(114, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(114, 16).Value = 1143
$ws.Cells.Item(114, 17).Value = 25

$ws.Cells.Item(115, 4).Value = 44924
$ws.Cells.Item(115, 9).Value = 'Primera'
$ws.Cells.Item(115, 10).Value = 205
$ws.Cells.Item(115, 11).Value = 43000
$ws.Cells.Item(115, 12).Value = 45000
$ws.Cells.Item(115, 13).Value = 44220
$ws.Cells.Item(115, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(115, 15).Value = 'Región del Maule'
$ws.Cells.Item(115, 16).Value = 1769
$ws.Cells.Item(115, 17).Value = 25

$ws.Cells.Item(116, 4).Value = 44973
$ws.Cells.Item(116, 9).Value = 'Primera'
$ws.Cells.Item(116, 10).Value = 80
$ws.Cells.Item(116, 11).Value = 30000
$ws.Cells.Item(116, 12).Value = 30000
$ws.Cells.Item(116, 13).Value = 30000
$ws.Cells.Item(116, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(116, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(116, 16).Value = 1200
$ws.Cells.Item(116, 17).Value = 25

$ws.Cells.Item(117, 4).Value = 44245
$ws.Cells.Item(117, 9).Value = 'Primera'
$ws.Cells.Item(117, 10).Value = 235
$ws.Cells.Item(117, 11).Value = 28000
$ws.Cells.Item(117, 12).Value = 29000
$ws.Cells.Item(117, 13).Value = 28532
$ws.Cells.Item(117, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(117, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(117, 16).Value = 1141
$ws.Cells.Item(117, 17).Value = 25

$ws.Cells.Item(118, 4).Value = 44558
$ws.Cells.Item(118, 9).Value = 'Primera'
$ws.Cells.Item(118, 10).Value = 28
$ws.Cells.Item(118, 11).Value = 30000
$ws.Cells.Item(118, 12).Value = 30000
$ws.Cells.Item(118, 13).Value = 30000
$ws.Cells.Item(118, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(118, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(118, 16).Value = 1200
$ws.Cells.Item(118, 17).Value = 25

$ws.Cells.Item(119, 4).Value = 44910
$ws.Cells.Item(119, 9).Value = 'Primera'
$ws.Cells.Item(119, 10).Value = 15
$ws.Cells.Item(119, 11).Value = 50000
$ws.Cells.Item(119, 12).Value = 50000
$ws.Cells.Item(119, 13).Value = 50000
$ws.Cells.Item(119, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(119, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(119, 16).Value = 2000
$ws.Cells.Item(119, 17).Value = 25

$ws.Cells.Item(120, 4).Value = 44980
$ws.Cells.Item(120, 9).Value = 'Primera'
$ws.Cells.Item(120, 10).Value = 45
$ws.Cells.Item(120, 11).Value = 33000
$ws.Cells.Item(120, 12).Value = 33000
$ws.Cells.Item(120, 13).Value = 33000
$ws.Cells.Item(120, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(120, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(120, 16).Value = 1320
$ws.Cells.Item(120, 17).Value = 25

$ws.Cells.Item(121, 4).Value = 44232
$ws.Cells.Item(121, 9).Value = 'Primera'
$ws.Cells.Item(121, 10).Value = 110
$ws.Cells.Item(121, 11).Value = 30000
$ws.Cells.Item(121, 12).Value = 30000
$ws.Cells.Item(121, 13).Value = 30000
$ws.Cells.Item(121, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(121, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(121, 16).Value = 1200
$ws.Cells.Item(121, 17).Value = 25

$ws.Cells.Item(122, 4).Value = 44613
$ws.Cells.Item(122, 9).Value = 'Primera'
$ws.Cells.Item(122, 10).Value = 70
$ws.Cells.Item(122, 11).Value = 25000
$ws.Cells.Item(122, 12).Value = 28000
$ws.Cells.Item(122, 13).Value = 26286
$ws.Cells.Item(122, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(122, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(122, 16).Value = 1051
$ws.Cells.Item(122, 17).Value = 25
